$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the taxReference value for the row-2 record (column I = "taxReference")
$ws.Range("I2").Value = "A555557"

# Update the active selection to match the edited cell area (I3, as in the target workbook)
$ws.Range("I3").Select()
